$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9
$ws.Range("G9").Value = 1.47
$ws.Range("I9").Value = 5.4
$ws.Range("N9").Value = 1.47
$ws.Range("O9").Value = 2.32
$ws.Range("R9").Value = 1.6
$ws.Range("S9").Value = 2.05
$ws.Range("T9").Value = 9.25
$ws.Range("U9").Value = 8.5
$ws.Range("W9").Value = 11
$ws.Range("X9").Value = 10.75
$ws.Range("Y9").Value = 20
$ws.Range("Z9").Value = 16.5
$ws.Range("AB9").Value = 15
$ws.Range("AC9").Value = 55
$ws.Range("AD9").Value = 19.5
$ws.Range("AE9").Value = 37
$ws.Range("AF9").Value = 17
$ws.Range("AG9").Value = 100
$ws.Range("AH9").Value = 45
$ws.Range("AI9").Value = 40
$ws.Range("AJ9").Value = 350

# Row 10
$ws.Range("G10").Value = 1.88
$ws.Range("H10").Value = 3.85
$ws.Range("I10").Value = 3.3
$ws.Range("N10").Value = 1.47
$ws.Range("O10").Value = 2.32
$ws.Range("S10").Value = 2.35
$ws.Range("T10").Value = 11
$ws.Range("U10").Value = 11.5
$ws.Range("W10").Value = 17.5
$ws.Range("X10").Value = 13.5
$ws.Range("Y10").Value = 18.5
$ws.Range("Z10").Value = 17.5
$ws.Range("AA10").Value = 8.25
$ws.Range("AB10").Value = 12
$ws.Range("AC10").Value = 37
$ws.Range("AD10").Value = 15
$ws.Range("AE10").Value = 22
$ws.Range("AF10").Value = 11.75
$ws.Range("AH10").Value = 25

# Row 12
$ws.Range("H12").Value = 2.92
$ws.Range("K12").Value = 6.2
$ws.Range("M12").Value = 2.8
$ws.Range("N12").Value = 2.12
$ws.Range("O12").Value = 1.65
$ws.Range("W12").Value = 20
$ws.Range("Z12").Value = 6.2
$ws.Range("AA12").Value = 5.8
$ws.Range("AG12").Value = 55

# Row 17
$ws.Range("K17").Value = 8
$ws.Range("N17").Value = 2.35
$ws.Range("O17").Value = 1.57
$ws.Range("AJ17").Value = 1000
